$d = $word.ActiveDocument

# Locate the paragraph ending with the "Il soupire..." stage direction so the
# insertion point is robust even if paragraph indices shift.
$target = $null
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Il soupire et se dirige vers le Ma*tre*") {
        $target = $candidate
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '(Il soupire...)' paragraph"
}

# Insert a new blank BodyText paragraph right after it.
[void]$target.Range.InsertParagraphAfter()

# Insert a second new blank BodyText paragraph right after the one we just
# added; this is the paragraph that will become "FOND NOIR".
$blankIndex = $targetIndex + 1
$blank = $d.Paragraphs.Item($blankIndex)
[void]$blank.Range.InsertParagraphAfter()

# The freshly created paragraph (now directly after $blank) gets replaced
# with properly bold-formatted "FOND NOIR" text via a raw OOXML fragment so
# both the run and the paragraph mark carry <w:b/><w:bCs/>.
$foundNoirIndex = $blankIndex + 1
$foundNoirPara = $d.Paragraphs.Item($foundNoirIndex)
$insertionPoint = $foundNoirPara.Range.Duplicate
$insertionPoint.Collapse(1)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>FOND NOIR</w:t></w:r></w:p>'
[void]$insertionPoint.InsertXML($xml)
